# Apply updated dSF (column F) values for a handful of rows.
# These reflect a "repull/push all data" refresh where the dSF figures
# were recalculated (mean calculation) and differ from the dS0 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = -4
    6  = -7
    9  = 1
    11 = -10
    14 = -2
    15 = -14
    18 = -9
    20 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
